$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Result" header in column D
$ws.Range("D1").Value = "Result"

# The computed results were typed/pasted in as text (carrying a trailing
# tab + newline), so recreate that exactly: build each value with a
# formula, then convert it to a literal value via copy / paste-values so
# it lands as plain text rather than being auto-coerced to a number.
$results = [ordered]@{
    2 = "8"
    3 = "10"
    4 = "18"
    5 = "0"
    6 = "5"
    7 = "4"
}

foreach ($r in $results.Keys) {
    $cell = $ws.Range("D$r")
    $cell.Formula = "=""" + $results[$r] + """&CHAR(9)&CHAR(10)"
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}
$excel.CutCopyMode = $false

# Extend the Result column formatting (wrap text) a few rows past the data
$ws.Range("D2:D15").WrapText = $true

# Make sure the trailing, otherwise-empty cells in column D actually exist
# in the sheet (so they get written out, still carrying the wrap style)
for ($r = 8; $r -le 15; $r++) {
    $ws.Range("D$r").Value = $null
}

# Give the data rows extra height to accommodate the wrapped text
for ($r = 2; $r -le 7; $r++) {
    $ws.Rows.Item($r).RowHeight = 30
}
